$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 175.6326516151428
$ws.Range("C2").Value = 2.783552277093058
$ws.Range("D2").Value = 1.637707805633545
$ws.Range("E2").Value = 0.1229069258522101

$ws.Range("B3").Value = 344.0073203086853
$ws.Range("C3").Value = 2.153519106192344
$ws.Range("D3").Value = 1.665207242965698
$ws.Range("E3").Value = 0.1956941583614969

$ws.Range("B4").Value = 689.87491106987
$ws.Range("C4").Value = 5.338644944978371
$ws.Range("D4").Value = 1.633451986312866
$ws.Range("E4").Value = 0.114246595484618

$ws.Range("B5").Value = 178.3612462520599
$ws.Range("C5").Value = 2.837976637256868
$ws.Range("D5").Value = 1.604144334793091
$ws.Range("E5").Value = 0.04503669638468374

$ws.Range("B6").Value = 347.7904633998871
$ws.Range("C6").Value = 1.99740067499088
$ws.Range("D6").Value = 1.62903528213501
$ws.Range("E6").Value = 0.09638317368159781

$ws.Range("B7").Value = 694.2723443508148
$ws.Range("C7").Value = 1.242448364147366
$ws.Range("D7").Value = 1.713821458816528
$ws.Range("E7").Value = 0.205582875125334

$ws.Range("B8").Value = 178.9650043487549
$ws.Range("C8").Value = 0.7040634891888485
$ws.Range("D8").Value = 1.682955503463745
$ws.Range("E8").Value = 0.1436807621405558

$ws.Range("B9").Value = 348.6202109336853
$ws.Range("C9").Value = 1.907522674659164
$ws.Range("D9").Value = 1.732587862014771
$ws.Range("E9").Value = 0.222965759323674

$ws.Range("B10").Value = 690.2973899841309
$ws.Range("C10").Value = 2.431504339423691
$ws.Range("D10").Value = 1.660965967178345
$ws.Range("E10").Value = 0.1926655467842077

$ws.Range("B11").Value = 292.2304166793823
$ws.Range("C11").Value = 1.637775234940798
$ws.Range("D11").Value = 1.60712251663208
$ws.Range("E11").Value = 0.07266480794412744

$ws.Range("B12").Value = 573.4947679519653
$ws.Range("C12").Value = 1.987377647349666
$ws.Range("D12").Value = 1.638145923614502
$ws.Range("E12").Value = 0.2057294217775834

$ws.Range("B13").Value = 1136.419561052322
$ws.Range("C13").Value = 4.600802567336125
$ws.Range("D13").Value = 1.608305406570435
$ws.Range("E13").Value = 0.1077966539801735

$ws.Range("B14").Value = 290.8895393371582
$ws.Range("C14").Value = 1.39301218110744
$ws.Range("D14").Value = 1.544513845443725
$ws.Range("E14").Value = 0.0504706129109295

$ws.Range("B15").Value = 577.1466466426849
$ws.Range("C15").Value = 2.645893417815718
$ws.Range("D15").Value = 1.633072376251221
$ws.Range("E15").Value = 0.1465522197038826

$ws.Range("B16").Value = 1158.428579950333
$ws.Range("C16").Value = 4.820982864573175
$ws.Range("D16").Value = 1.705912733078003
$ws.Range("E16").Value = 0.187090630444331

$ws.Range("B17").Value = 298.0820672988892
$ws.Range("C17").Value = 2.21276743162678
$ws.Range("D17").Value = 1.62284426689148
$ws.Range("E17").Value = 0.1323213834214018

$ws.Range("B18").Value = 580.9934951782227
$ws.Range("C18").Value = 4.707513736109631
$ws.Range("D18").Value = 1.55162181854248
$ws.Range("E18").Value = 0.04567344913507335

$ws.Range("B19").Value = 935.7997517585754
$ws.Range("C19").Value = 34.40619817686384
$ws.Range("D19").Value = 1.022534561157227
$ws.Range("E19").Value = 0.2185300749728398
